$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Gender values were stored with title-case text ("Female"/"Male"); update them
# to lower-case ("female"/"male") to match the values now being inserted into
# the database (see commit message: "wrote insert code from df to db").
$ws.Range("H2").Value = "female"
$ws.Range("H3").Value = "female"
$ws.Range("H4").Value = "male"

# Reflect the new selection left behind in the saved file.
$ws.Range("H5").Select()
